## Regenerate merged AHB files
## - Rename the "_old"/"_new" comparison-column headers to "_FV2210"/"_FV2304"
## - Turn the sheet's data range into an Excel Table ("Table1")
## - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Columns A:J were "<Name>_old", L:U were "<Name>_new"; column K is "diff".
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    # columns 1-10 => A-J ("_old" -> "_FV2210")
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2210"
    # columns 12-21 => L-U ("_new" -> "_FV2304")
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2304"
}

# --- 2. Convert the used range into a native Excel Table -------------------
$tableRange = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
